$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G4's existing value (time portion changed)
$ws.Range("G4").Value = 42606.498668981483

# New row 5
$ws.Range("A5").Value = 9820.09
$ws.Range("B5").Value = 9875.39
$ws.Range("C5").Value = 316.81
$ws.Range("D5").Value = 318.58999999999997
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = 0.56000000000000005
$ws.Range("G5").Value = 42606.585590277777
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
$ws.Range("H5").Value = $false

# New row 6
$ws.Range("A6").Value = 9752.33
$ws.Range("B6").Value = 9820.09
$ws.Range("C6").Value = 316.81
$ws.Range("D6").Value = 319
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = 0.69
$ws.Range("G6").Value = 42606.586712962962
$ws.Range("G6").NumberFormat = "m/d/yy h:mm"
$ws.Range("H6").Value = $false

# New row 7
$ws.Range("A7").Value = 9697.7199999999993
$ws.Range("B7").Value = 9752.33
$ws.Range("C7").Value = 316.81
$ws.Range("D7").Value = 318.58999999999997
$ws.Range("E7").Value = $true
$ws.Range("F7").Value = 0.56000000000000005
$ws.Range("G7").Value = 42606.587824074071
$ws.Range("G7").NumberFormat = "m/d/yy h:mm"
$ws.Range("H7").Value = $false
